# add teleport after finish normal level
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (BattleScene) value label: "UI/Ice2Prefab" -> "Ice2Prefab"
$ws.Range("C3:C5").Value = "Ice2Prefab"

# Update column D (BossLocation) value label: "UI/Ice3Prefab" -> "Image/UI/Ice3"
$ws.Range("D3:D5").Value = "Image/UI/Ice3"

# New header row (row 1) for the added columns
$ws.Range("E1").Value = "NormalEnemyAmount"
$ws.Range("F1").Value = "BossEnemyAmount"
$ws.Range("G1").Value = "SpawnOne"
$ws.Range("H1").Value = "SpawnTwo"
$ws.Range("I1").Value = "SpawnThree"

# New sub-header row (row 2) for the added columns
$ws.Range("E2").Value = "Number of enemy"
$ws.Range("F2").Value = "Number of enemy"
$ws.Range("G2").Value = "Spawn Enemy"
$ws.Range("H2").Value = "Spawn Enemy"
$ws.Range("I2").Value = "Spawn Enemy"

# New data rows 3-5 for the added columns
$ws.Range("E3:E5").Value = 3
$ws.Range("F3:F5").Value = 1
$ws.Range("G3:G5").Value = "Model/cheche(main)"
$ws.Range("H3:H5").Value = "Model/cheche(main)"
$ws.Range("I3:I5").Value = "Model/cheche(main)"

# Selection moved to I10 as recorded in the saved file
$ws.Range("I10").Select()
